$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# Helper: turn a MatchScorecard URL into a bare match-code string,
# forcing the destination cell to remain TEXT (not auto-coerced to a
# Number) by assigning through a leading single-quote, exactly like
# typing '1234 into Excel.
# ------------------------------------------------------------------
function Set-TextValue($cell, $text) {
    $cell.Value = "'" + $text
}

# ==================================================================
# 1. "ODI Batting" sheet (existing sheet #1): MATCH_CARD_LINK ->
#    MATCH_CODE, drop the (empty) INNING_NUMBER placeholder cells,
#    and collapse each URL down to its trailing match code.
# ==================================================================
$batting = $wb.Worksheets.Item("ODI Batting")

$batting.Range("D1").Value = "MATCH_CODE"

$lastRow = 50
for ($r = 2; $r -le $lastRow; $r++) {
    $linkCell = $batting.Cells.Item($r, 4)
    $link = $linkCell.Value2
    if ($link -ne $null) {
        $marker = "MatchCode="
        $idx = $link.IndexOf($marker)
        if ($idx -ge 0) {
            $code = $link.Substring($idx + $marker.Length)
            Set-TextValue $linkCell $code
        }
    }

    $inningCell = $batting.Cells.Item($r, 2)
    if ($inningCell.Value2 -eq $null) {
        $inningCell.ClearContents()
    }
}

# ==================================================================
# 2. "ODI Bowling" sheet (existing sheet #2): same MATCH_CARD_LINK ->
#    MATCH_CODE treatment on column B.
# ==================================================================
$bowling = $wb.Worksheets.Item("ODI Bowling")

$bowling.Range("B1").Value = "MATCH_CODE"

for ($r = 2; $r -le $lastRow; $r++) {
    $linkCell = $bowling.Cells.Item($r, 2)
    $link = $linkCell.Value2
    if ($link -ne $null) {
        $marker = "MatchCode="
        $idx = $link.IndexOf($marker)
        if ($idx -ge 0) {
            $code = $link.Substring($idx + $marker.Length)
            Set-TextValue $linkCell $code
        }
    }
}

# ==================================================================
# 3. New "Player Info" sheet, inserted before "ODI Batting" so the
#    final tab order is: Player Info, ODI Batting, ODI Bowling,
#    ODI Batting Extra.
# ==================================================================
$playerInfo = $wb.Worksheets.Add($batting)
$playerInfo.Name = "Player Info"

$headerRange = $playerInfo.Range("A1:D1")
$headerRange.Font.Bold = $true
$headerRange.Borders.LineStyle = 1
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160

$playerInfo.Range("A1").Value = "ID"
$playerInfo.Range("B1").Value = "NAME"
$playerInfo.Range("C1").Value = "BATTING_HAND"
$playerInfo.Range("D1").Value = "BOWL_STYLE"

Set-TextValue $playerInfo.Range("A2") "3914"
$playerInfo.Range("B2").Value = "Aththachchi Nuwan Pradeep Roshan Fernando"
$playerInfo.Range("C2").Value = "Right Handed"
$playerInfo.Range("D2").Value = "Right Arm Fast Medium"

# ==================================================================
# 4. New "ODI Batting Extra" sheet, appended after "ODI Bowling".
#    Re-fetch "ODI Bowling" by name since inserting "Player Info"
#    above shifted every sheet's numeric position - a stale
#    positional handle would now point at the wrong tab.
# ==================================================================
$bowlingNow = $wb.Worksheets.Item("ODI Bowling")
$extra = $wb.Worksheets.Add($null, $bowlingNow)
$extra.Name = "ODI Batting Extra"

$extraHeader = $extra.Range("A1:F1")
$extraHeader.Font.Bold = $true
$extraHeader.Borders.LineStyle = 1
$extraHeader.HorizontalAlignment = -4108
$extraHeader.VerticalAlignment = -4160

$extra.Range("A1").Value = "MATCH_CODE"
$extra.Range("B1").Value = "BATTING_POSITION"
$extra.Range("C1").Value = "NUM_4"
$extra.Range("D1").Value = "NUM_6"
$extra.Range("E1").Value = "PERCENT_RUNS_OF_TOTAL"
$extra.Range("F1").Value = "MAN_OF_MATCH"

# MATCH_CODE, BATTING_POSITION, NUM_4, NUM_6, PERCENT_RUNS_OF_TOTAL, MAN_OF_MATCH
$extraRows = @(
    @("4210", 11,   $null, $null, $null,    "NO"),
    @("4211", $null,$null, $null, $null,    "NO"),
    @("4231", $null,$null, $null, $null,    "NO"),
    @("4232", 11,   "0",   "0",   $null,    "YES"),
    @("4233", 11,   "0",   "0",   "0.43%",  "NO"),
    @("4302", $null,$null, $null, $null,    "NO"),
    @("4309", $null,$null, $null, $null,    "NO"),
    @("4322", 10,   "0",   "0",   "0.34%",  "NO"),
    @("4331", $null,$null, $null, $null,    "NO"),
    @("4356", 11,   "0",   "0",   $null,    "NO"),
    @("4357", 11,   "0",   "0",   "0.43%",  "NO"),
    @("4375", $null,$null, $null, $null,    "NO"),
    @("4376", $null,$null, $null, $null,    "NO"),
    @("4413", $null,$null, $null, $null,    $null),
    @("4414", $null,$null, $null, $null,    $null),
    @("4417", $null,$null, $null, $null,    $null),
    @("4449", $null,$null, $null, $null,    $null),
    @("4450", $null,$null, $null, $null,    $null),
    @("4521", $null,$null, $null, $null,    $null),
    @("4523", $null,$null, $null, $null,    $null)
)

$row = 2
foreach ($data in $extraRows) {
    Set-TextValue $extra.Cells.Item($row, 1) $data[0]

    if ($data[1] -ne $null) {
        $extra.Cells.Item($row, 2).Value = $data[1]
    }
    if ($data[2] -ne $null) {
        Set-TextValue $extra.Cells.Item($row, 3) $data[2]
    }
    if ($data[3] -ne $null) {
        Set-TextValue $extra.Cells.Item($row, 4) $data[3]
    }
    if ($data[4] -ne $null) {
        Set-TextValue $extra.Cells.Item($row, 5) $data[4]
    }
    if ($data[5] -ne $null) {
        $extra.Cells.Item($row, 6).Value = $data[5]
    }

    $row = $row + 1
}

Write-Output "done"
